$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Wins" (AD1), "Losses" (AE1), "Ties" (AF1) ---
# Copy the existing header formatting (from the adjacent header cell AC1,
# style s="1": bold font, thin border, centered/top aligned) onto the new
# header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-55: season record for every player's team ---
# Every row in this sheet gets the same team season record: 74 wins,
# 88 losses, 0 ties.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 74  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 88  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
